$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.875.39"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "3.105.90"
$ws.Range("E3").Value = "  -0.38%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "577.25"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "173.77"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "3.108.33"
$ws.Range("E8").Value = "  -0.25%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.515"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.30%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.39"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.47%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.152"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.73%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.475"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.12%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000241"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -3.15%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "35.95"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -3.33%  "
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "3.618.75"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "66.750.22"
$ws.Range("E17").Value = "  -0.41%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "17.08"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +4.29%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.00"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("D20").Value = "3.100.03"
$ws.Range("E20").Value = "  -0.59%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "485.37"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.13%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.82"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.49%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.692"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.01%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "83.52"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.72%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "12.72"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -3.76%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.24"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.10%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.11"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.62%  "
$ws.Range("E28").Value = "  +0.00%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.03"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.60%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.27"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.66%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.59"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.06%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "28.10"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.65%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.112"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("D34").Value = "0.0₃0939"
$ws.Range("E34").Value = "  -0.71%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.19%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "48.08"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.39%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.61"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.57%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.946"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.86%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "49.13"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.02%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.310"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.16%  "
$ws.Range("E41").Value = "  -0.46%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.97"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -4.02%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "8.30"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.56%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.67"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.61%  "
$ws.Range("D45").Value = "2.803.56"
$ws.Range("E45").Value = "  -0.69%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.0347"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.57%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "370.11"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.33%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "134.38"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.88%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "24.49"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.03%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.20"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.59%  "
